$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the new "time (in ms)" column (F) mirroring column E ---
$ws.Range("F5").Value = "time (in ms)"
$ws.Range("F6:F12").Value = 0

# Copy formatting (border + left alignment + general style) from column E onto column F
$ws.Range("E5:E12").Copy()
$ws.Range("F5:F12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Update the chart ---
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart

# Name the existing series ("Sans index")
$ser1 = $chart.SeriesCollection().Item(1)
$ser1.Name = "Sans index"

# Add the second series ("Avec index") plotting the new column F
$ser2 = $chart.SeriesCollection().NewSeries()
$ser2.Formula = "=SERIES(`"Avec index`",,Feuil1!`$F`$6:`$F`$12,2)"
$ser2.Border.Weight = 2.25
$ser2.Border.Color = 3243501   # matches theme accent2 (ED7D31) line color

# Show a legend on the right so the two series can be told apart
$chart.HasLegend = $true
$chart.Legend.Position = -4152  # xlLegendPositionRight

# --- 3. Enlarge the chart to make room for the legend ---
$chartObj.Width = 447.7
$chartObj.Height = 279

# --- 4. Leave the selection where the author left it ---
$ws.Range("O22").Select()
